$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from the last existing
# header cell (AC1) onto the three new header cells so they match the rest
# of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (same constant values for every player row 2-47)
$ws.Range("AD2:AD47").Value = 89
$ws.Range("AE2:AE47").Value = 73
$ws.Range("AF2:AF47").Value = 0
